$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 899; everything from 899 downward shifts to 900.
$ws.Rows.Item(899).Insert()

# Populate the newly inserted row 899 with the new record.
$ws.Cells.Item(899, 1).Value = 4
$ws.Cells.Item(899, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(899, 3).Value = "Los Lagos"
$ws.Cells.Item(899, 4).Value = 45223
$ws.Cells.Item(899, 5).Value = 10
$ws.Cells.Item(899, 6).Value = "Fruta"
$ws.Cells.Item(899, 7).Value = 100102
$ws.Cells.Item(899, 8).Value = "Cítricos"
$ws.Cells.Item(899, 9).Value = 100102003
$ws.Cells.Item(899, 10).Value = "Limón"
$ws.Cells.Item(899, 11).Value = "Sin especificar"
$ws.Cells.Item(899, 12).Value = "1a plateado"
$ws.Cells.Item(899, 13).Value = 1200
$ws.Cells.Item(899, 14).Value = 14000
$ws.Cells.Item(899, 15).Value = 15000
$ws.Cells.Item(899, 16).Value = 14500
$ws.Cells.Item(899, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(899, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(899, 19).Value = 806
$ws.Cells.Item(899, 20).Value = 18
